$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet ("fix xlsx data" generic rename to Sheet1)
$ws.Name = "Sheet1"

# Fix the surface-tension (ST) target band: 30-38 -> 32-38
$ws.Range("G3").Value = "32－38"

# Row 14 (DP_test_03): replace the evaluation comment and priority mark
$ws.Range("I14").Value = "画像濃度、表面張力が予測値から上振れすれば目標値達成の可能性がある`n上振れの程度もRMSEよりも小さい量であり、ある程度可能性は残されていると推察できる"
$ws.Range("J14").Value = "△"

# Row 15 (DP_test_04): priority mark upgraded
$ws.Range("J15").Value = "◎"

# Row 18 (DP_test_07): evaluation comment now also calls out 耐擦過性, priority mark downgraded
$ws.Range("I18").Value = "画像濃度,耐擦過性が予測値から上振れすれば目標値達成の可能性がある`n上振れの程度もRMSEよりも小さい量であり、ある程度可能性は残されていると推察できる"
$ws.Range("J18").Value = "△"

# Cosmetic sheet-view tweaks captured in the diff
$ws.Range("I1").EntireColumn.ColumnWidth = 63.9140625

$window = $excel.ActiveWindow
$window.DisplayGridlines = $false
